$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163
$xlPasteFormats = -4122

# 1) Free up the merged cells in the region we are about to restructure
#    (data rows 7-12, totals row 13, footer row 14) so the row-by-row
#    copy below can write into every individual cell.
$ws.Range("A7:N14").UnMerge()

# 2) Shift rows 7..14 down to 8..15 (bottom-up so nothing is clobbered
#    before it is read): row 14 (footer) -> 15, row 13 (totals) -> 14,
#    and the 6 medicine rows 12..7 -> 13..8.
for ($src = 14; $src -ge 7; $src--) {
    $dst = $src + 1
    $srcRange = $ws.Range("A" + $src + ":N" + $src)
    $dstRange = $ws.Range("A" + $dst + ":N" + $dst)

    $srcRange.Copy()
    $dstRange.PasteSpecial($xlPasteValues)

    $srcRange.Copy()
    $dstRange.PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = 0

# 3) Write the new medicine row into row 7 (between HAIRVOGINE and
#    METFORMIN, alphabetically) -- formatting/style already matches the
#    other data rows because row 7 previously held METFORMIN before the
#    shift above.
$ws.Cells.Item(7, 1).Value2 = 4
$ws.Cells.Item(7, 2).Value2 = "LYSE 0.65% NASAL DROPS 15 ML"
$ws.Cells.Item(7, 8).Value2 = "0:0"
$ws.Cells.Item(7, 12).Value2 = 30
$ws.Cells.Item(7, 14).Value2 = "2:0"

# 4) A column just holds the running position 1..10 and is unaffected by
#    the shift, but the brand-new last data row (13) needs it set.
$ws.Cells.Item(13, 1).Value2 = 10

# 5) Update the grand total (now on row 14) to include the new row.
$ws.Cells.Item(14, 11).Value2 = 609

# 6) Re-merge the data-row cell groups for rows 7-13.
for ($r = 7; $r -le 13; $r++) {
    $ws.Range("B" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()
}

# 7) Re-merge the totals row (now 14) and footer row (now 15).
$ws.Range("K14:N14").Merge()
$ws.Range("A15:E15").Merge()
$ws.Range("F15:G15").Merge()
$ws.Range("I15:N15").Merge()

# 8) Row heights: the new data row 13 keeps the standard 25.5pt height,
#    and the footer (now row 15) shrinks slightly to 16.5pt as it does
#    whenever Excel reflows this report after inserting a row.
$ws.Rows(13).RowHeight = 25.5
$ws.Rows(14).RowHeight = 25.5
$ws.Rows(15).RowHeight = 16.5
